$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
Write-Output $lastSheet.Name
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Notes"
Write-Output $wb.Worksheets.Count
foreach ($w in $wb.Worksheets) {
    Write-Output $w.Name
}
